# Applies the 5/12 literature-review update to Sheet 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# --- New "notes" (column D) for existing rows, plus three new rows at the
# bottom of the list. Order matches the author's original edit/typing order
# (this determines shared-string append order on save).
$ws.Range("D33").Value = "bayesian model of CE in value-based choice"
$ws.Range("D38").Value = "Context Dependence and Aggregation in Disaggregate Choice Analysis"
$ws.Range("D36").Value = "range-normalization model"
$ws.Range("D35").Value = "CE choice model"
$ws.Range("D40").Value = "mlba"
$ws.Range("D41").Value = "review - Theories of context effects in multialternative, multiattribute choice"
$ws.Range("D43").Value = "mlca"
$ws.Range("D44").Value = "2n ary choice tree"
$ws.Range("D45").Value = "Stochastic choice: An optimizing neuroeconomic model"
$ws.Range("D87").Value = "Testing the effect of time pressure on asymmetric dominance and compromise decoys in choice"

$ws.Range("A151").Value = "Pettibone 2000"
$ws.Range("D151").Value = "Examining {Models} of {Nondominated} {Decoy} {Effects} across {Judgment} and {Choice}"

$ws.Range("D93").Value = "Impact of choice set complexity on decoy effects, also eye-tracking and comparison "

$ws.Range("A152").Value = "Wedell & Pettibone 1996"
$ws.Range("D152").Value = "Using Judgments to Understand Decoy Effects in Choice"

$ws.Range("A153").Value = "Windschitl & Chambers 2004"
$ws.Range("D153").Value = "The Dud-Alternative Effect in Likelihood Judgment"

# --- View state: scrolled/resized window + selection moved to new last row -
$ws.Application.ActiveWindow.ScrollRow = 135
$ws.Range("E153").Select()
